$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 48, shifting rows 48:63 down to 49:64
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new weekly record
$ws.Cells.Item(48, 1).Value = 4
$ws.Cells.Item(48, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(48, 3).Value = "Los Lagos"
$ws.Cells.Item(48, 4).Value = 44917
$ws.Cells.Item(48, 5).Value = 10
$ws.Cells.Item(48, 6).Value = "Fruta"
$ws.Cells.Item(48, 7).Value = 100103
$ws.Cells.Item(48, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(48, 9).Value = 100103003
$ws.Cells.Item(48, 10).Value = "Damasco"
$ws.Cells.Item(48, 11).Value = "Castle Brite"
$ws.Cells.Item(48, 12).Value = "Primera"
$ws.Cells.Item(48, 13).Value = 300
$ws.Cells.Item(48, 14).Value = 19000
$ws.Cells.Item(48, 15).Value = 20000
$ws.Cells.Item(48, 16).Value = 19500
$ws.Cells.Item(48, 17).Value = "`$/caja 16 kilos"
$ws.Cells.Item(48, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(48, 19).Value = 1219
$ws.Cells.Item(48, 20).Value = 16
